# Auto-generated script applying scheduled-runner market-data refresh
# to the Leve profit calculator sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each touched row: H/I/J are refreshed raw market prices; K=I*F, L=J*F,
# M=E-K (only if I<>0), N=-2*E-L (only if J<>0) are the recomputed leve profit figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1258.75
$ws.Range("I9").Value = 160.36363
$ws.Range("J9").Value = 3675.2
$ws.Range("K9").Value = 160.36363
$ws.Range("L9").Value = 3675.2
$ws.Range("M9").Value = 8.63637
$ws.Range("N9").Value = -4013.2
$ws.Range("H92").Value = 883.36664
$ws.Range("I92").Value = 1065.4
$ws.Range("J92").Value = 519.3
$ws.Range("K92").Value = 1065.4
$ws.Range("L92").Value = 519.3
$ws.Range("M92").Value = 182.5999999999999
$ws.Range("N92").Value = -3015.3
$ws.Range("H98").Value = 2909.0908
$ws.Range("I98").Value = 2935.0527
$ws.Range("J98").Value = 2744.6667
$ws.Range("K98").Value = 2935.0527
$ws.Range("L98").Value = 2744.6667
$ws.Range("M98").Value = -1437.0527000000002
$ws.Range("N98").Value = -5740.6667
$ws.Range("H122").Value = 2909.0908
$ws.Range("I122").Value = 2935.0527
$ws.Range("J122").Value = 2744.6667
$ws.Range("K122").Value = 8805.1581
$ws.Range("L122").Value = 8234.000100000001
$ws.Range("M122").Value = -6355.158100000001
$ws.Range("N122").Value = -13134.000100000001
$ws.Range("H132").Value = 4832.0713
$ws.Range("I132").Value = 4834.923
$ws.Range("J132").Value = 4795
$ws.Range("K132").Value = 14504.769
$ws.Range("L132").Value = 14385
$ws.Range("M132").Value = -11974.769
$ws.Range("N132").Value = -19445
$ws.Range("H137").Value = 2280
$ws.Range("I137").Value = 1147.8334
$ws.Range("J137").Value = 2802.5386
$ws.Range("K137").Value = 3443.5002
$ws.Range("L137").Value = 8407.6158
$ws.Range("M137").Value = -893.5002
$ws.Range("N137").Value = -13507.6158
$ws.Range("H138").Value = 3430.1453
$ws.Range("I138").Value = 1009.8
$ws.Range("J138").Value = 3642.456
$ws.Range("K138").Value = 3029.3999999999996
$ws.Range("L138").Value = 10927.368
$ws.Range("M138").Value = 2110.6000000000004
$ws.Range("N138").Value = -21207.368000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 12614
$ws.Range("I31").Value = 12614
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 12614
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -12320
$ws.Range("N31").ClearContents()
$ws.Range("H32").Value = 3297096
$ws.Range("I32").Value = 3455378.5
$ws.Range("J32").Value = 1002000
$ws.Range("K32").Value = 3455378.5
$ws.Range("L32").Value = 1002000
$ws.Range("M32").Value = -3455091.5
$ws.Range("N32").Value = -1002574
$ws.Range("H74").Value = 2560.3333
$ws.Range("I74").Value = 2334.2856
$ws.Range("J74").Value = 2876.8
$ws.Range("K74").Value = 2334.2856
$ws.Range("L74").Value = 2876.8
$ws.Range("M74").Value = -1460.2856000000002
$ws.Range("N74").Value = -4624.8
$ws.Range("H77").Value = 2560.3333
$ws.Range("I77").Value = 2334.2856
$ws.Range("J77").Value = 2876.8
$ws.Range("K77").Value = 11671.428
$ws.Range("L77").Value = 14384
$ws.Range("M77").Value = -7303.428
$ws.Range("N77").Value = -23120
$ws.Range("H104").Value = 45577.6
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 45577.6
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 45577.6
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -52565.6
$ws.Range("H132").Value = 3792.8
$ws.Range("I132").Value = 3883.5
$ws.Range("J132").Value = 3732.3333
$ws.Range("K132").Value = 11650.5
$ws.Range("L132").Value = 11196.999899999999
$ws.Range("M132").Value = -9120.5
$ws.Range("N132").Value = -16256.999899999999
$ws.Range("H139").Value = 10806.667
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 10806.667
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 10806.667
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -21086.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7628.5
$ws.Range("I94").Value = 3748.8572
$ws.Range("J94").Value = 14417.875
$ws.Range("K94").Value = 3748.8572
$ws.Range("L94").Value = 14417.875
$ws.Range("M94").Value = -3297.8572
$ws.Range("N94").Value = -15319.875
$ws.Range("H107").Value = 5506
$ws.Range("I107").Value = 4670.3335
$ws.Range("J107").Value = 8013
$ws.Range("K107").Value = 4670.3335
$ws.Range("L107").Value = 8013
$ws.Range("M107").Value = -2750.3334999999997
$ws.Range("N107").Value = -11853
$ws.Range("H134").Value = 2962.375
$ws.Range("I134").Value = 3180.3215
$ws.Range("J134").Value = 1436.75
$ws.Range("K134").Value = 9540.9645
$ws.Range("L134").Value = 4310.25
$ws.Range("M134").Value = -7005.9645
$ws.Range("N134").Value = -9380.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1514.7778
$ws.Range("I31").Value = 1184.7073
$ws.Range("J31").Value = 4898
$ws.Range("K31").Value = 1184.7073
$ws.Range("L31").Value = 4898
$ws.Range("M31").Value = -889.7073
$ws.Range("N31").Value = -5488
$ws.Range("H34").Value = 1514.7778
$ws.Range("I34").Value = 1184.7073
$ws.Range("J34").Value = 4898
$ws.Range("K34").Value = 1184.7073
$ws.Range("L34").Value = 4898
$ws.Range("M34").Value = -982.7073
$ws.Range("N34").Value = -5302
$ws.Range("H47").Value = 38971
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 38971
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 38971
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -40103
$ws.Range("H58").Value = 4011.45
$ws.Range("I58").Value = 3256.0908
$ws.Range("J58").Value = 4934.6665
$ws.Range("K58").Value = 3256.0908
$ws.Range("L58").Value = 4934.6665
$ws.Range("M58").Value = -3053.0908
$ws.Range("N58").Value = -5340.6665
$ws.Range("H107").Value = 2377.9644
$ws.Range("I107").Value = 1969.4736
$ws.Range("J107").Value = 3240.3333
$ws.Range("K107").Value = 1969.4736
$ws.Range("L107").Value = 3240.3333
$ws.Range("M107").Value = -49.47360000000003
$ws.Range("N107").Value = -7080.3333
$ws.Range("H132").Value = 1072.5
$ws.Range("I132").Value = 1245
$ws.Range("J132").Value = 555
$ws.Range("K132").Value = 3735
$ws.Range("L132").Value = 1665
$ws.Range("M132").Value = -1205
$ws.Range("N132").Value = -6725
$ws.Range("H136").Value = 4011.45
$ws.Range("I136").Value = 3256.0908
$ws.Range("J136").Value = 4934.6665
$ws.Range("K136").Value = 9768.2724
$ws.Range("L136").Value = 14803.999500000002
$ws.Range("M136").Value = -7218.2724
$ws.Range("N136").Value = -19903.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3707.3333
$ws.Range("I55").Value = 2990
$ws.Range("J55").Value = 3772.5454
$ws.Range("K55").Value = 8970
$ws.Range("L55").Value = 11317.6362
$ws.Range("M55").Value = -8793
$ws.Range("N55").Value = -11671.6362
$ws.Range("H60").Value = 4571.25
$ws.Range("I60").Value = 5093
$ws.Range("J60").Value = 3006
$ws.Range("K60").Value = 15279
$ws.Range("L60").Value = 9018
$ws.Range("M60").Value = -15028
$ws.Range("N60").Value = -9520
$ws.Range("H131").Value = 1680.0488
$ws.Range("I131").Value = 1126.5714
$ws.Range("J131").Value = 1794
$ws.Range("K131").Value = 3379.7142000000003
$ws.Range("L131").Value = 5382
$ws.Range("M131").Value = 1660.2857999999997
$ws.Range("N131").Value = -15462

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6816.2354
$ws.Range("I70").Value = 6534.385
$ws.Range("J70").Value = 7732.25
$ws.Range("K70").Value = 6534.385
$ws.Range("L70").Value = 7732.25
$ws.Range("M70").Value = -6264.385
$ws.Range("N70").Value = -8272.25
$ws.Range("H73").Value = 6816.2354
$ws.Range("I73").Value = 6534.385
$ws.Range("J73").Value = 7732.25
$ws.Range("K73").Value = 6534.385
$ws.Range("L73").Value = 7732.25
$ws.Range("M73").Value = -5598.385
$ws.Range("N73").Value = -9604.25
$ws.Range("H80").Value = 1966.25
$ws.Range("I80").Value = 1938.5
$ws.Range("J80").Value = 1994
$ws.Range("K80").Value = 1938.5
$ws.Range("L80").Value = 1994
$ws.Range("M80").Value = -940.5
$ws.Range("N80").Value = -3990
$ws.Range("H83").Value = 1966.25
$ws.Range("I83").Value = 1938.5
$ws.Range("J83").Value = 1994
$ws.Range("K83").Value = 9692.5
$ws.Range("L83").Value = 9970
$ws.Range("M83").Value = -4700.5
$ws.Range("N83").Value = -19954
$ws.Range("H97").Value = 3103
$ws.Range("I97").Value = 2799.9524
$ws.Range("J97").Value = 3898.5
$ws.Range("K97").Value = 2799.9524
$ws.Range("L97").Value = 3898.5
$ws.Range("M97").Value = -2303.9524
$ws.Range("N97").Value = -4890.5
$ws.Range("H122").Value = 4758.1816
$ws.Range("I122").Value = 9135.333
$ws.Range("J122").Value = 3116.75
$ws.Range("K122").Value = 27405.999000000003
$ws.Range("L122").Value = 9350.25
$ws.Range("M122").Value = -24955.999000000003
$ws.Range("N122").Value = -14250.25
$ws.Range("H132").Value = 1770.8572
$ws.Range("I132").Value = 1599.25
$ws.Range("J132").Value = 1999.6666
$ws.Range("K132").Value = 4797.75
$ws.Range("L132").Value = 5998.9998
$ws.Range("M132").Value = -2267.75
$ws.Range("N132").Value = -11058.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 38999.5
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 38999.5
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 38999.5
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -39223.5
$ws.Range("H22").Value = 16825.176
$ws.Range("I22").Value = 1809.2727
$ws.Range("J22").Value = 44354.332
$ws.Range("K22").Value = 1809.2727
$ws.Range("L22").Value = 44354.332
$ws.Range("M22").Value = -1514.2727
$ws.Range("N22").Value = -44944.332
$ws.Range("H27").Value = 16825.176
$ws.Range("I27").Value = 1809.2727
$ws.Range("J27").Value = 44354.332
$ws.Range("K27").Value = 1809.2727
$ws.Range("L27").Value = 44354.332
$ws.Range("M27").Value = -1702.2727
$ws.Range("N27").Value = -44568.332
$ws.Range("H68").Value = 3217.5715
$ws.Range("I68").Value = 3829.1
$ws.Range("J68").Value = 1688.75
$ws.Range("K68").Value = 3829.1
$ws.Range("L68").Value = 1688.75
$ws.Range("M68").Value = -3080.1
$ws.Range("N68").Value = -3186.75
$ws.Range("H71").Value = 3217.5715
$ws.Range("I71").Value = 3829.1
$ws.Range("J71").Value = 1688.75
$ws.Range("K71").Value = 19145.5
$ws.Range("L71").Value = 8443.75
$ws.Range("M71").Value = -15401.5
$ws.Range("N71").Value = -15931.75
$ws.Range("H122").Value = 4521.5557
$ws.Range("I122").Value = 2075
$ws.Range("J122").Value = 6478.8
$ws.Range("K122").Value = 6225
$ws.Range("L122").Value = 19436.4
$ws.Range("M122").Value = -3775
$ws.Range("N122").Value = -24336.4
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10500
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7970
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 23706.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 23706.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 23706.5
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -25578.5
$ws.Range("H77").Value = 23706.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 23706.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 71119.5
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -80479.5
$ws.Range("H136").Value = 787.7436
$ws.Range("I136").Value = 609.1111
$ws.Range("J136").Value = 2931.3333
$ws.Range("K136").Value = 1827.3332999999998
$ws.Range("L136").Value = 8793.999899999999
$ws.Range("M136").Value = 722.6667000000002
$ws.Range("N136").Value = -13893.999899999999
